$d = $word.ActiveDocument

# The four inline logo pictures live in the document's two headers and two
# footers. Their wp:docPr / pic:cNvPr "name" attributes were swapped:
#   - the Pearson logo pictures (descr ...PearsonLogo.png) go from
#     "image2.png" to "image1.png"
#   - the BTEC logo pictures (descr "BTec_Logo-Orange") go from
#     "image1.jpg" to "image2.jpg"
# InlineShape does not expose a settable Name property in the Word object
# model, so round-trip the canonical WordOpenXML for the whole document
# (which includes headers/footers) and patch the two name attributes
# textually, then write it back.

$xml = $d.Content.WordOpenXML

$xml = $xml.Replace('name="image1.jpg"', 'name="image2.jpg"')
$xml = $xml.Replace('name="image2.png"', 'name="image1.png"')

$d.Content.WordOpenXML = $xml
